# Updates the "cryptos" price/volume snapshot (Price in column D, Volume(1h)
# in column E) to a newer refresh, and fixes three rows whose Coin/Link pairs
# had been listed in the wrong relative order (rows 20/21, 37/38, 44/45).
#
# Column D holds prices formatted as look-alike-numeric text (e.g. "574.08",
# "0.387", even "1.00") stored as plain strings in the source workbook (no
# thousands grouping, just literal text). A naive `.Value = "574.08"` would
# be auto-coerced by Excel into the *number* 574.08 (and would silently
# collapse something like "1.00" into the number 1). Stamping the range as
# Text ("@") before the assignment -- then putting the style back to Normal
# so no stray per-cell formatting is left behind -- keeps these as text,
# matching the workbook's original representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "61.004.48"
$ws.Range("E2").Value = "  -1.75%  "
Set-TextValue "D3" "3.392.28"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "574.08"
$ws.Range("E5").Value = "  -0.69%  "
Set-TextValue "D6" "137.62"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.04%  "
Set-TextValue "D8" "3.392.65"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  -1.26%  "
Set-TextValue "D10" "7.65"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  -2.08%  "
Set-TextValue "D12" "0.387"
$ws.Range("E12").Value = "  -1.61%  "
Set-TextValue "D13" "3.969.92"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  +0.50%  "
Set-TextValue "D15" "26.18"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("E16").Value = "  -2.70%  "
Set-TextValue "D17" "3.388.30"
$ws.Range("E17").Value = "  -0.46%  "
Set-TextValue "D18" "61.149.26"
$ws.Range("E18").Value = "  -1.49%  "
Set-TextValue "D19" "14.03"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "9.49"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D21" "5.82"
$ws.Range("E21").Value = "  -1.13%  "
Set-TextValue "D22" "377.00"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("E23").Value = "  -2.72%  "
Set-TextValue "D24" "3.527.23"
$ws.Range("E24").Value = "  -0.46%  "
Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -1.78%  "
Set-TextValue "D27" "71.18"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  +11.93%  "
Set-TextValue "D29" "0.175"
$ws.Range("E29").Value = "  +9.22%  "
Set-TextValue "D30" "7.53"
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("E34").Value = "  +0.01%  "
Set-TextValue "D35" "23.75"
$ws.Range("E35").Value = "  +0.71%  "
Set-TextValue "D36" "5.22"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D37" "6.87"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "1.55"
$ws.Range("E38").Value = "  -0.67%  "
Set-TextValue "D39" "164.38"
$ws.Range("E39").Value = "  +0.06%  "
Set-TextValue "D40" "0.0762"
$ws.Range("E40").Value = "  -3.43%  "
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  +0.03%  "
Set-TextValue "D42" "0.780"
$ws.Range("E42").Value = "  -0.69%  "
Set-TextValue "D43" "1.71"
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D44" "41.68"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D45" "4.42"
$ws.Range("E45").Value = "  -1.28%  "
Set-TextValue "D46" "1.20"
$ws.Range("E46").Value = "  -2.68%  "
Set-TextValue "D47" "24.51"
$ws.Range("E47").Value = "  -2.64%  "
Set-TextValue "D48" "2.472.46"
$ws.Range("E48").Value = "  +4.29%  "
Set-TextValue "D49" "23.19"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  -2.68%  "
Set-TextValue "D51" "2.45"
$ws.Range("E51").Value = "  +6.26%  "
